$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Status now "ok" (Cliente interface text unchanged)
$ws.Range("C3").Value = "ok"

# Rows 5-7: text unchanged, just touch so shared strings are reused in order
$ws.Range("B5").Value = "Construir  Interfaces dos servicos Vacinacao ( Valter)"
$ws.Range("B6").Value = "Construir  Interfaces dos servicos Exames  (Jaime)"
$ws.Range("B7").Value = "Construir  Interfaces dos servicos cirurgia  ( Valter)"

# Row 8: historico -> append (Jaime)
$ws.Range("B8").Value = "Construir  Interfaces do historico ( Jaime)"

# Row 9: relatorios -> append (Jaime e Valter)
$ws.Range("B9").Value = "Fazer relatorios ( Jaime e Valter)"

# Row 10: login Senha -> append (Jaime)
$ws.Range("B10").Value = "Tela de login (Senha)  (Jaime)"

# Row 11: new text - Tela de Loading (Valter)
$ws.Range("B11").Value = "Tela de Loading (Valter)"

# Row 12: new row - insercao/actualizacao/eliminacao (Jaime e Valter)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Fazer a insercao, actualizacao e eliminacao  de varios dados  (Jaime e Valter)"

# Row 4: Animal interface text updated + now marked "ok" (written last so its
# shared string is appended after the other new strings)
$ws.Range("B4").Value = "Construir  Interfaces de cadastro e consulta Animal (Jaime ) "
$ws.Range("C4").Value = "ok"

# Update selection to match target (B6)
$ws.Range("B6").Select()
